$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: change status from "offen" to "optional" ---
# Copy format from a cell that already uses the "optional" style (B7), then set text.
$ws.Range("B7").Copy()
$ws.Range("B20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B20").Value = "optional"

# --- Row 21: change status from "offen" to "done" ---
# Copy format from a cell that already uses the "done" style (B2), then set text.
$ws.Range("B2").Copy()
$ws.Range("B21").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B21").Value = "done"

# --- New Row 22: new TODO item (keeps "offen" style, like B19) ---
$ws.Range("A22").Value = "E-Mail -> Texte überarbeiten (eventuell mit HTML was machen?)"
$ws.Range("B19").Copy()
$ws.Range("B22").PasteSpecial(-4122) # xlPasteFormats (copy "offen" style)
$ws.Range("B22").Value = "offen"

# --- Row 19: add note in column D ---
$ws.Range("D19").Value = "Tests"

# Clear clipboard/marquee state left over from the copy/paste operations
$excel.CutCopyMode = $false

# Update selection to match target workbook state
$ws.Range("D19").Select()
